$wb = $excel.ActiveWorkbook

# The "Spain" sheet is structurally the closest existing template for the
# three new market sheets: same row layout/styles/merges/margins, the
# ht="28.8" rows 3-5, and column widths/bestFit flags on columns A/B that
# already match the target new sheets exactly.

# --- Russia -----------------------------------------------------------
$wb.Worksheets.Item("Spain").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$russia = $wb.Worksheets.Item($wb.Worksheets.Count)
$russia.Name = "Russia"
$russia.Columns.Item(4).ColumnWidth = 7.6
$russia.Range("B4").Value = "NGC-2929/T2900"
$russia.Range("B2").Value = "Russia Market"

# --- Finland ------------------------------------------------------------
$wb.Worksheets.Item("Spain").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$finland = $wb.Worksheets.Item($wb.Worksheets.Count)
$finland.Name = "Finland"
$finland.Columns.Item(4).ColumnWidth = 7.6
$finland.Range("B4").Value = "NGC-3130/T2943"
$finland.Range("B2").Value = "Finland Market"

# --- Hungary --------------------------------------------------------------
$wb.Worksheets.Item("Spain").Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$hungary = $wb.Worksheets.Item($wb.Worksheets.Count)
$hungary.Name = "Hungary"
$hungary.Columns.Item(4).ColumnWidth = 7.6
$hungary.Range("B4").Value = "NGC-3104/T2992"
$hungary.Range("B2").Value = "Hungary Market"

# Russia & Finland keep the plain A1:D10 selection (not the active tab).
$russia.Activate()
$russia.Range("A1:D10").Select()
$finland.Activate()
$finland.Range("A1:D10").Select()

# Hungary ends up as the active sheet/tab, cursor parked at H15.
$hungary.Activate()
$hungary.Range("H15").Select()
